# Update "Max Potential Capacity by Source" workbook per latest data refresh.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")
$wsMpcbs = $wb.Worksheets.Item("MPCbS")

# --- About sheet: clarify methodology note ---
$wsAbout.Range("A11").Value = "We used the technical potential of renewable energy sources"

# --- Data sheet: refreshed max potential capacity figures ---
$wsData.Range("B5").Value = 12000
$wsData.Range("B6").Value = 352000
$wsData.Range("B7").Value = 2409000
$wsData.Range("B9").Value = 10000
$wsData.Range("B10").Value = 1259000
$wsData.Range("B14").Value = 387000

# --- Selections / active sheet to match the saved view state ---
$wsData.Range("B17").Select() | Out-Null
$wsMpcbs.Range("B5").Select() | Out-Null
$wsAbout.Range("B10").Select() | Out-Null
$wsAbout.Activate() | Out-Null
